$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 715, pushing existing rows 715-768 down to 717-770.
$ws.Rows("715:716").Insert()

# New row 715: Primera quality, week of 2022-08-10 (serial 44783)
$ws.Range("A715").Value = 3
$ws.Range("B715").Value = "Femacal de La Calera"
$ws.Range("C715").Value = "Coquimbo"
$ws.Range("D715").Value = 44783
$ws.Range("E715").Value = 5
$ws.Range("F715").Value = 100112023
$ws.Range("G715").Value = "Brócoli"
$ws.Range("H715").Value = "Sin especificar"
$ws.Range("I715").Value = "Primera"
$ws.Range("J715").Value = 1730
$ws.Range("K715").Value = 850
$ws.Range("L715").Value = 900
$ws.Range("M715").Value = 875
$ws.Range("N715").Value = "`$/unidad"
$ws.Range("O715").Value = "Provincia de Quillota"
$ws.Range("P715").Value = 875
$ws.Range("Q715").Value = 1
$ws.Range("R715").Value = "Hortaliza"

# New row 716: Segunda quality, same week (serial 44783)
$ws.Range("A716").Value = 3
$ws.Range("B716").Value = "Femacal de La Calera"
$ws.Range("C716").Value = "Coquimbo"
$ws.Range("D716").Value = 44783
$ws.Range("E716").Value = 5
$ws.Range("F716").Value = 100112023
$ws.Range("G716").Value = "Brócoli"
$ws.Range("H716").Value = "Sin especificar"
$ws.Range("I716").Value = "Segunda"
$ws.Range("J716").Value = 1750
$ws.Range("K716").Value = 700
$ws.Range("L716").Value = 750
$ws.Range("M716").Value = 724
$ws.Range("N716").Value = "`$/unidad"
$ws.Range("O716").Value = "Provincia de Quillota"
$ws.Range("P716").Value = 724
$ws.Range("Q716").Value = 1
$ws.Range("R716").Value = "Hortaliza"
